$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 34328.332
$ws.Range("J128").Value = 34328.332
$ws.Range("L128").Value = 34328.332
$ws.Range("N128").Value = -44288.332
$ws.Range("H130").Value = 37850
$ws.Range("J130").Value = 37850
$ws.Range("L130").Value = 37850
$ws.Range("N130").Value = -47890
$ws.Range("H132").Value = 2326.842
$ws.Range("I132").Value = 2013.7174
$ws.Range("J132").Value = 3636.2727
$ws.Range("K132").Value = 6041.1522
$ws.Range("L132").Value = 10908.8181
$ws.Range("M132").Value = -3511.1522
$ws.Range("N132").Value = -15968.8181
$ws.Range("H136").Value = 88483.164
$ws.Range("J136").Value = 88483.164
$ws.Range("L136").Value = 88483.164
$ws.Range("N136").Value = -98683.164
$ws.Range("H139").Value = 46790
$ws.Range("J139").Value = 46790
$ws.Range("L139").Value = 46790
$ws.Range("N139").Value = -57070
$ws.Range("H140").Value = 64624.547
$ws.Range("I140").Value = 135000
$ws.Range("J140").Value = 57587
$ws.Range("K140").Value = 135000
$ws.Range("L140").Value = 57587
$ws.Range("M140").Value = -129820
$ws.Range("N140").Value = -67947

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22676.973
$ws.Range("I32").Value = 4441.113
$ws.Range("J32").Value = 125460.91
$ws.Range("K32").Value = 4441.113
$ws.Range("L32").Value = 125460.91
$ws.Range("M32").Value = -4154.113
$ws.Range("N32").Value = -126034.91
$ws.Range("H121").Value = 13167.917
$ws.Range("J121").Value = 13167.917
$ws.Range("L121").Value = 13167.917
$ws.Range("N121").Value = -16661.917
$ws.Range("H122").Value = 2219.5334
$ws.Range("I122").Value = 2163
$ws.Range("J122").Value = 2375
$ws.Range("K122").Value = 6489
$ws.Range("L122").Value = 7125
$ws.Range("M122").Value = -4039
$ws.Range("N122").Value = -12025
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("H132").Value = 3372.1785
$ws.Range("I132").Value = 1895.7222
$ws.Range("J132").Value = 6029.8
$ws.Range("K132").Value = 5687.1666
$ws.Range("L132").Value = 18089.4
$ws.Range("M132").Value = -3157.1666
$ws.Range("N132").Value = -23149.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1336.7241
$ws.Range("I107").Value = 1011.875
$ws.Range("J107").Value = 2896
$ws.Range("K107").Value = 1011.875
$ws.Range("L107").Value = 2896
$ws.Range("M107").Value = 908.125
$ws.Range("N107").Value = -6736
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("H130").Value = 35540
$ws.Range("J130").Value = 35540
$ws.Range("L130").Value = 35540
$ws.Range("N130").Value = -45580

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2118.7964
$ws.Range("I31").Value = 1421.2963
$ws.Range("J31").Value = 2816.2964
$ws.Range("K31").Value = 1421.2963
$ws.Range("L31").Value = 2816.2964
$ws.Range("M31").Value = -1126.2963
$ws.Range("N31").Value = -3406.2964
$ws.Range("H34").Value = 2118.7964
$ws.Range("I34").Value = 1421.2963
$ws.Range("J34").Value = 2816.2964
$ws.Range("K34").Value = 1421.2963
$ws.Range("L34").Value = 2816.2964
$ws.Range("M34").Value = -1219.2963
$ws.Range("N34").Value = -3220.2964
$ws.Range("H99").Value = 1631.8462
$ws.Range("J99").Value = 1833.3334
$ws.Range("L99").Value = 1833.3334
$ws.Range("N99").Value = -4829.3334
$ws.Range("H126").Value = 1631.8462
$ws.Range("J126").Value = 1833.3334
$ws.Range("L126").Value = 5500.0002
$ws.Range("N126").Value = -10440.0002
$ws.Range("H132").Value = 3145.2
$ws.Range("I132").Value = 2905.3333
$ws.Range("J132").Value = 4104.6665
$ws.Range("K132").Value = 8715.999899999999
$ws.Range("L132").Value = 12313.9995
$ws.Range("M132").Value = -6185.999899999999
$ws.Range("N132").Value = -17373.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 73.083336
$ws.Range("I38").Value = 68
$ws.Range("J38").Value = 80.2
$ws.Range("K38").Value = 204
$ws.Range("L38").Value = 240.6
$ws.Range("M38").Value = 143
$ws.Range("N38").Value = -934.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1738.6285
$ws.Range("I102").Value = 1264.8
$ws.Range("J102").Value = 2923.2
$ws.Range("K102").Value = 1264.8
$ws.Range("L102").Value = 2923.2
$ws.Range("M102").Value = 357.2
$ws.Range("N102").Value = -6167.2
$ws.Range("H122").Value = 1409.95
$ws.Range("I122").Value = 1352.1538
$ws.Range("J122").Value = 1517.2858
$ws.Range("K122").Value = 4056.4614
$ws.Range("L122").Value = 4551.857400000001
$ws.Range("M122").Value = -1606.4614
$ws.Range("N122").Value = -9451.857400000001
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("H128").Value = 49742.5
$ws.Range("J128").Value = 49742.5
$ws.Range("L128").Value = 49742.5
$ws.Range("N128").Value = -59702.5
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040
$ws.Range("H132").Value = 3090.3215
$ws.Range("I132").Value = 2647.2307
$ws.Range("J132").Value = 3474.3333
$ws.Range("K132").Value = 7941.6921
$ws.Range("L132").Value = 10422.9999
$ws.Range("M132").Value = -5411.6921
$ws.Range("N132").Value = -15482.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14287543
$ws.Range("I7").Value = 16668133
$ws.Range("K7").Value = 16668133
$ws.Range("M7").Value = -16668021
$ws.Range("H40").Value = 1614.4445
$ws.Range("I40").Value = 1439.3334
$ws.Range("J40").Value = 2490
$ws.Range("K40").Value = 1439.3334
$ws.Range("L40").Value = 2490
$ws.Range("M40").Value = -1303.3334
$ws.Range("N40").Value = -2762
$ws.Range("H61").Value = 2238.0833
$ws.Range("I61").Value = 1493.1428
$ws.Range("J61").Value = 3281
$ws.Range("K61").Value = 1493.1428
$ws.Range("L61").Value = 3281
$ws.Range("M61").Value = -1291.1428
$ws.Range("N61").Value = -3685
$ws.Range("H95").Value = 29929.334
$ws.Range("J95").Value = 29929.334
$ws.Range("L95").Value = 29929.334
$ws.Range("N95").Value = -35421.334
$ws.Range("H113").Value = 2238.0833
$ws.Range("I113").Value = 1493.1428
$ws.Range("J113").Value = 3281
$ws.Range("K113").Value = 1493.1428
$ws.Range("L113").Value = 3281
$ws.Range("M113").Value = 676.8571999999999
$ws.Range("N113").Value = -7621
$ws.Range("H118").Value = 33117
$ws.Range("J118").Value = 33117
$ws.Range("L118").Value = 33117
$ws.Range("N118").Value = -36431
$ws.Range("H126").Value = 14287543
$ws.Range("I126").Value = 16668133
$ws.Range("K126").Value = 50004399
$ws.Range("M126").Value = -50001929

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 40585.875
$ws.Range("J105").Value = 40585.875
$ws.Range("L105").Value = 40585.875
$ws.Range("N105").Value = -47573.875
$ws.Range("H113").Value = 56544.277
$ws.Range("I113").Value = 91229.17999999999
$ws.Range("J113").Value = 2039.4286
$ws.Range("K113").Value = 273687.54
$ws.Range("L113").Value = 6118.2858
$ws.Range("M113").Value = -271517.54
$ws.Range("N113").Value = -10458.2858
$ws.Range("H122").Value = 10068.16
$ws.Range("I122").Value = 11714.2
$ws.Range("J122").Value = 3484
$ws.Range("K122").Value = 35142.60000000001
$ws.Range("L122").Value = 10452
$ws.Range("M122").Value = -32692.60000000001
$ws.Range("N122").Value = -15352
$ws.Range("H126").Value = 2703.3125
$ws.Range("I126").Value = 1930.3
$ws.Range("J126").Value = 3991.6667
$ws.Range("K126").Value = 5790.9
$ws.Range("L126").Value = 11975.0001
$ws.Range("M126").Value = -3320.9
$ws.Range("N126").Value = -16915.0001
$ws.Range("H132").Value = 1313.0483
$ws.Range("I132").Value = 1222.849
$ws.Range("J132").Value = 1844.2222
$ws.Range("K132").Value = 3668.547
$ws.Range("L132").Value = 5532.6666
$ws.Range("M132").Value = -1138.547
$ws.Range("N132").Value = -10592.6666
